# Trade #5 closed at 2026-02-17 07:57:47 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly-closed trade #5 to both
# the "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a cell as literal text, never letting the engine's
# auto-date-detection turn date-looking strings (e.g. "2026-02-17")
# into a date serial number / date-formatted cell.
# ---------------------------------------------------------------------
function Set-TextCell($range, $text) {
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value2 = 1199.94   # Current Capital
$summary.Range("B4").Value2 = -0.06     # Total P&L $
$summary.Range("B5").Value2 = -0.24     # Total P&L %
$summary.Range("B6").Value2 = 5         # Total Trades
$summary.Range("B8").Value2 = 4         # Losing Trades
$summary.Range("B9").Value2 = 20        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value2 = 99.94      # Capital
$status.Range("D4").Value2 = 5          # Trades
$status.Range("E4").Value2 = -0.06      # P&L $
$status.Range("F4").Value2 = -0.06      # P&L %
$status.Range("G4").Value2 = 20         # Win Rate %

# ---------------------------------------------------------------------
# New trade #5 row, appended to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A6").Value2 = 5
    Set-TextCell $ws.Range("B6") "2026-02-17"
    $ws.Range("C6").Value2 = "07:57:41"
    $ws.Range("D6").Value2 = "MarketMaking"
    $ws.Range("E6").Value2 = "DOWN"
    $ws.Range("F6").Value2 = 0.79
    $ws.Range("G6").Value2 = 0.78
    $ws.Range("H6").Value2 = "CLOSED"
    $ws.Range("I6").Value2 = -1.2658
    $ws.Range("J6").Value2 = -0.01
    $ws.Range("K6").Value2 = 99.94
    $ws.Range("L6").Value2 = 0
    $ws.Range("M6").Value2 = 0
    $ws.Range("N6").Value2 = 0.6
    $ws.Range("O6").Value2 = "Normal spread capture: 19600 bps"
    $ws.Range("P6").Value2 = "early_exit"
    $ws.Range("Q6").Value2 = 0.13
}
